$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (before account 004927044 / CINTIA) and
# populate it with the new account (EVANGELINA). The account number has a
# leading zero, so force text entry via a temporary Text format, then drop
# the format again so the cell matches its unstyled neighbours.
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "005646524"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "EVANGELINA"
$ws.Cells.Item(2, 3).Value = 1000000

# Remove the four rows that were dropped from the export. Match by the
# account number in column A and walk bottom-up so row indices stay valid
# as rows are deleted.
$accountsToRemove = @("004419141", "004420763", "008004799", "004384258")

$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
for ($r = $lastRow; $r -ge 2; $r--) {
    $acct = $ws.Cells.Item($r, 1).Value2
    if ($accountsToRemove -contains $acct) {
        $ws.Rows.Item($r).Delete()
    }
}
